$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.340.75"
$ws.Range("E2").Value = "  +1.64%  "

$ws.Range("D3").Value = "3.896.90"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'528.58"
$ws.Range("E5").Value = "  +9.54%  "

$ws.Range("D6").Value = "'144.61"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("E7").Value = "  -1.59%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -3.27%  "

$ws.Range("E10").Value = "  -4.89%  "

$ws.Range("D11").Value = "'0.0000335"
$ws.Range("E11").Value = "  -5.59%  "

$ws.Range("D12").Value = "'42.04"
$ws.Range("E12").Value = "  -2.23%  "

$ws.Range("D13").Value = "4.519.10"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").Value = "'10.23"
$ws.Range("E14").Value = "  -2.51%  "

$ws.Range("D15").Value = "3.921.79"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").Value = "'13.97"
$ws.Range("E16").Value = "  -1.66%  "

$ws.Range("D17").Value = "'0.134"
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("D18").Value = "'1.21"
$ws.Range("E18").Value = "  +6.98%  "

$ws.Range("D19").Value = "'19.78"
$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").Value = "69.218.62"
$ws.Range("E20").Value = "  +1.43%  "

$ws.Range("D21").Value = "'425.49"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("E22").Value = "  -6.39%  "

$ws.Range("D23").Value = "'88.12"
$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("D24").Value = "'14.10"
$ws.Range("E24").Value = "  -4.60%  "

$ws.Range("E25").Value = "  +10.00%  "

$ws.Range("D26").Value = "'11.39"
$ws.Range("E26").Value = "  -8.11%  "

$ws.Range("D27").Value = "'10.57"
$ws.Range("E27").Value = "  -4.03%  "

$ws.Range("D28").Value = "'36.34"
$ws.Range("E28").Value = "  -2.41%  "

$ws.Range("D29").Value = "'686.17"
$ws.Range("E29").Value = "  -4.50%  "

$ws.Range("D30").Value = "'13.13"
$ws.Range("E30").Value = "  -2.63%  "

$ws.Range("E31").Value = "  -3.16%  "

$ws.Range("D32").Value = "'2.83"
$ws.Range("E32").Value = "  -2.86%  "

$ws.Range("D33").Value = "'68.75"
$ws.Range("E33").Value = "  +11.21%  "

$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0884"
$ws.Range("E34").Value = "  +1.43%  "

$ws.Range("B35").Value = "TheGraph"
$ws.Range("C35").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D35").Value = "'0.434"
$ws.Range("E35").Value = "  +8.65%  "

$ws.Range("E36").Value = "  -2.05%  "

$ws.Range("D37").Value = "'39.89"
$ws.Range("E37").Value = "  -2.40%  "

$ws.Range("E38").Value = "  +1.86%  "

$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").Value = "'3.29"
$ws.Range("E40").Value = "  +7.17%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").Value = "'0.0480"
$ws.Range("E42").Value = "  -3.64%  "

$ws.Range("D43").Value = "'3.20"
$ws.Range("E43").Value = "  +8.42%  "

$ws.Range("E44").Value = "  -7.05%  "

$ws.Range("D45").Value = "'3.40"
$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("D46").Value = "'0.000286"
$ws.Range("E46").Value = "  +17.72%  "

$ws.Range("E47").Value = "  -1.80%  "

$ws.Range("D48").Value = "'3.00"
$ws.Range("E48").Value = "  +6.71%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'146.01"
$ws.Range("E49").Value = "  +0.98%  "

$ws.Range("D50").Value = "2.742.39"
$ws.Range("E50").Value = "  +14.61%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0344"
$ws.Range("E51").Value = "  -4.09%  "

